# Extracted a common code block into a method.
#
# Adds a literal value in B1, which feeds the existing B3 formula
# (=A2+B1), and adds a new summary formula in A5 that totals the
# original A1:B3 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New input cell used by the existing B3 formula (=A2+B1).
$ws.Range("B1").Value = 13

# New row with a formula that sums the whole original A1:B3 block.
$ws.Range("A5").Formula = "=SUM(A1:B3)"

# Match the author's final selection/active cell.
$ws.Range("A5").Select()
